$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.867.26'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '1.825.21'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6910'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07608'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3014'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.34'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07713'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.73%  '
$ws.Range('D12').Value = '1.829.33'
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.033'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6704'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.352'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008264'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').Value = '28.861.26'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.35%  '
$ws.Range('D20').Value = '2.077.31'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.374'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1468'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.697'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.528'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.177'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.115'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.188'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05082'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7445'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.805'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.137'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.682'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01834'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = '1.197.53'
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.664'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9147'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9991'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5161'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.976.89'
$ws.Range('E45').Value = '  -3.05%  '
$ws.Range('E46').Value = '  -5.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.437'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.182'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -13.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.718'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '62.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -13.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4179'
$ws.Range('D51').Style = 'Normal'
